# Apply cryptocurrency price/volume updates for Sun Jan 22 08:29:04 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$text) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = $origStyle
}

Set-TextValue "D2" "302.00"
Set-TextValue "D3" "37.53"
Set-TextValue "E3" "6.95%"
Set-TextValue "D4" "5.002"
Set-TextValue "E4" "-3.22%"
Set-TextValue "D5" "0.07830"
Set-TextValue "E5" "0.53%"
Set-TextValue "D6" "2.192"
Set-TextValue "E6" "-7.67%"
Set-TextValue "D7" "8.034"
Set-TextValue "E7" "-0.30%"
Set-TextValue "D8" "4.018"
Set-TextValue "E8" "1.17%"
Set-TextValue "D9" "0.9102"
Set-TextValue "E9" "-1.43%"
Set-TextValue "D10" "0.09693"
Set-TextValue "E10" "-2.85%"
Set-TextValue "D11" "0.1890"
Set-TextValue "E11" "3.77%"
Set-TextValue "D12" "0.08538"
Set-TextValue "E12" "-1.33%"
Set-TextValue "D13" "0.03525"
Set-TextValue "E13" "6.23%"
Set-TextValue "D14" "0.09959"
Set-TextValue "E14" "0.48%"
Set-TextValue "D15" "0.001486"
Set-TextValue "E15" "-0.06%"
Set-TextValue "D16" "0.005725"
Set-TextValue "E16" "0.29%"
Set-TextValue "D17" "3.464"
Set-TextValue "E17" "-0.02%"
Set-TextValue "D18" "2.068"
Set-TextValue "E18" "-4.58%"
Set-TextValue "E19" "2.58%"
Set-TextValue "D20" "0.1294"
Set-TextValue "E20" "-2.40%"
Set-TextValue "D21" "4.765"
Set-TextValue "E21" "10.44%"
Set-TextValue "D22" "0.2207"
Set-TextValue "E22" "-7.38%"
Set-TextValue "D23" "0.04633"
Set-TextValue "E23" "1.36%"
Set-TextValue "E24" "1.06%"
Set-TextValue "D25" "0.004799"
Set-TextValue "E25" "7.88%"
Set-TextValue "E26" "-7.54%"
Set-TextValue "D27" "0.0004756"
Set-TextValue "E27" "28.68%"
Set-TextValue "D39" "0.01750"
Set-TextValue "E39" "-1.70%"
Set-TextValue "D40" "0.04711"
Set-TextValue "E40" "-1.71%"
Set-TextValue "D41" "0.008058"
Set-TextValue "E41" "4.51%"
Set-TextValue "D42" "0.1390"
Set-TextValue "E42" "-1.52%"
Set-TextValue "D43" "0.007674"
Set-TextValue "E43" "7.19%"
Set-TextValue "E44" "-3.35%"
Set-TextValue "D45" "0.009915"
Set-TextValue "E45" "3.91%"
Set-TextValue "D46" "0.00006058"
Set-TextValue "E46" "-0.77%"
Set-TextValue "E47" "0.21%"
Set-TextValue "E50" "0.21%"
Set-TextValue "D51" "0.0002002"
Set-TextValue "E51" "0.21%"
